# Actualizando proyecto antes de comenzar agregar scripts
# Adds 8 new "DEC_01xx" test rows (73-80) to the data pool, keeping the
# trailing blank-row block / footer rows intact (shifted down), and moves
# the active selection to reflect where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 fresh rows right above the old "TC / Verity1.0" footer row (76)
# so the rows below (old 76 -> TC/Verity1.0, old 77 -> credentials) shift
# down to 85/86, and the new rows 76:84 inherit the formatting (style) of
# the row above, same as Excel does on a manual row insert.
$ws.Rows("76:84").Insert()

# New DEC_ codes continuing the existing sequence.
$decCodes = @("DEC_0156","DEC_0157","DEC_0158","DEC_0159","DEC_0160","DEC_0161","DEC_0162","DEC_0163")

for ($i = 0; $i -lt $decCodes.Length; $i++) {
    $r = 73 + $i
    $ws.Range("A$r").Value = $decCodes[$i]
    $ws.Range("B$r").Value = "13712759-8"
    $ws.Range("C$r").Value = "Verity1.1"
    $ws.Range("D$r`:J$r").Value = "SIN_DATO"
}

# Reflect the user's final scroll position / selection in the saved view.
$ws.Application.Goto($ws.Range("A45")) | Out-Null
$ws.Range("E74").Select() | Out-Null
